$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# fix: correct value typo in classification table
# Correct value 2023 (non-existing in AAA-Objektartenkatalog) to 2013
$ws.Range("B2").Value = 2013

# Reapply the default (Normal) cell style so B2 no longer carries the
# now-superfluous explicit style that only existed to hold the old value
$ws.Range("B2").Style = "Normal"
